$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L4").Value = "Multivalued"
$ws.Range("L4").Font.Bold = $true

$ws.Range("L5").Value = "FALSE"
$ws.Range("L6").Value = "FALSE"
$ws.Range("L7").Value = "FALSE"

$ws.Range("L4:L7").Select()
